$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.383.92"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "'3.600.75"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'580.04"
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("D6").Value = "'190.21"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D8").Value = "'3.596.86"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.186"
$ws.Range("E10").Value = "  +4.26%  "
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "'56.15"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "'0.0000312"
$ws.Range("E13").Value = "  +7.76%  "
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "'4.177.56"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "'19.81"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'3.598.52"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").Value = "'70.305.64"
$ws.Range("D19").Value = "'12.64"
$ws.Range("E19").Value = "  -0.31%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "'491.93"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").Value = "'19.60"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("E24").Value = "  -9.49%  "
$ws.Range("D25").Value = "'96.49"
$ws.Range("E25").Value = "  +6.10%  "
$ws.Range("E26").Value = "  -2.23%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  -4.51%  "
$ws.Range("D28").Value = "'11.02"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'9.39"
$ws.Range("E29").Value = "  -3.09%  "
$ws.Range("D30").Value = "'32.27"
$ws.Range("E30").Value = "  -2.44%  "
$ws.Range("D31").Value = "'7.61"
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").Value = "'12.28"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("D33").Value = "'66.24"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "'577.41"
$ws.Range("E35").Value = "  -7.75%  "
$ws.Range("D36").Value = "'38.80"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").Value = "'0.0₃0815"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "  +5.95%  "
$ws.Range("D41").Value = "'3.25"
$ws.Range("E41").Value = "  +17.68%  "
$ws.Range("D42").Value = "'3.48"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("E43").Value = "  -6.32%  "
$ws.Range("D44").Value = "'3.225.15"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("E45").Value = "  -3.07%  "
$ws.Range("D46").Value = "'0.0447"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'9.65"
$ws.Range("E47").Value = "  +5.18%  "
$ws.Range("D48").Value = "'3.40"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "'3.21"
$ws.Range("E51").Value = "  -3.31%  "
